# Fixed custom plots (the way OF is used was very badly wrong @_@).
# Added some custom fields for useful quantities in NKF1 and NKF6, PSC.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new rows needed for the new "OF" based plot rows ---
# Row 7 (new): Roll Rate (OF) - inserted right after the existing "Des Roll Rate" row (row 6)
$ws.Rows.Item(7).Insert()
# Row 10 (new, after the first insert shifted things down): Pitch Rate (OF) -
# inserted right after the "Des Pitch Rate" row (originally row 8, now row 9)
$ws.Rows.Item(10).Insert()

# --- Rewrite the whole data block (rows 6-15) to match the corrected table ---

function Set-PlotRow($r, $plotNo, $row, $col, $axesLayout, $xAxisLabel, $yAxisLabel, $legendStyle, $yChannel, $scaleFactor, $labelOverride, $unitOverride) {
    $ws.Cells.Item($r, 1).Value = $plotNo
    $ws.Cells.Item($r, 2).Value = $row
    $ws.Cells.Item($r, 3).Value = $col
    $ws.Cells.Item($r, 4).Value = $axesLayout
    $ws.Cells.Item($r, 5).Value = $xAxisLabel
    $ws.Cells.Item($r, 6).Value = $yAxisLabel
    $ws.Cells.Item($r, 7).Value = $legendStyle
    $ws.Cells.Item($r, 11).Value = $yChannel
    $ws.Cells.Item($r, 15).Value = $scaleFactor
    $ws.Cells.Item($r, 18).Value = $labelOverride
    $ws.Cells.Item($r, 19).Value = $unitOverride
}

# Plot 1 - Roll Rate
Set-PlotRow 6  1 2 1 "S" "Time [ s ]" "Roll~Rate~[~deg/s~]" "Vertical" "RATE/RDes" 1    "Des Roll Rate"  "deg/s"
Set-PlotRow 7  1 2 1 "S" "Time [ s ]" "Roll~Rate~[~deg/s~]" "Vertical" "OF/bodyX"  57.7 "Roll Rate (OF)" "deg/s"
Set-PlotRow 8  1 2 1 "S" "Time [ s ]" "Roll~Rate~[~deg/s~]" "Vertical" "RATE/R"    1    "Roll Rate"      "deg/s"

# Plot 2 - Pitch Rate
Set-PlotRow 9  2 2 1 "S" "Time [ s ]" "Pitch~Rate~[~deg/s~]" "Vertical" "RATE/PDes" 1    "Des Pitch Rate"  "deg/s"
Set-PlotRow 10 2 2 1 "S" "Time [ s ]" "Pitch~Rate~[~deg/s~]" "Vertical" "OF/bodyY"  57.7 "Pitch Rate (OF)" "deg/s"
Set-PlotRow 11 2 2 1 "S" "Time [ s ]" "Pitch~Rate~[~deg/s~]" "Vertical" "RATE/P"    1    "Pitch Rate"      "deg/s"

# Plot 3 - Yaw Rate
Set-PlotRow 12 3 2 2 "S" "Time [ s ]" "Yaw~Rate~[~deg/s~]" "Vertical" "RATE/YDes" 1 "Des Yaw Rate" "deg/s"
Set-PlotRow 13 3 2 2 "S" "Time [ s ]" "Yaw~Rate~[~deg/s~]" "Vertical" "RATE/Y"    1 "Yaw Rate"     "deg/s"

# Plot 4 - Acceleration
Set-PlotRow 14 4 2 2 "S" "Time [ s ]" "Acceleration~[~m/s/s~]" "Vertical" "RATE/ADes" 1 "Des Acceleration" "m/s/s"
Set-PlotRow 15 4 2 2 "S" "Time [ s ]" "Acceleration~[~m/s/s~]" "Vertical" "RATE/A"    1 "Acceleration"     "m/s/s"

# --- Update view: scroll to column E, select I14 ---
$ws.Range("I14").Select()
$ws.Application.ActiveWindow.ScrollColumn = 5
